$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "52.361.92"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3
$ws.Range("D3").Value = "2.926.74"
$ws.Range("E3").Value = "  +4.54%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.91"
$ws.Range("E5").Value = "  -0.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.46"
$ws.Range("E6").Value = "  +3.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.628"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.15"
$ws.Range("E10").Value = "  +0.10%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("E11").Value = "  +3.01%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.136"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.14"
$ws.Range("E13").Value = "  +0.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.88"
$ws.Range("E14").Value = "  +1.25%  "

# Row 15
$ws.Range("D15").Value = "3.382.50"
$ws.Range("E15").Value = "  +4.21%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.936.56"
$ws.Range("E16").Value = "  +4.66%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "  +6.65%  "

# Row 18
$ws.Range("D18").Value = "52.357.97"
$ws.Range("E18").Value = "  +1.46%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.73"
$ws.Range("E19").Value = "  -0.79%  "

# Row 20
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.35"
$ws.Range("E20").Value = "  +5.25%  "

# Row 21
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.60"
$ws.Range("E21").Value = "  +8.93%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +1.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.10"
$ws.Range("E23").Value = "  +1.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.34"
$ws.Range("E24").Value = "  +1.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +1.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.92"
$ws.Range("E26").Value = "  +3.13%  "

# Row 28
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.58"
$ws.Range("E29").Value = "  +2.13%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.59"
$ws.Range("E30").Value = "  +3.41%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.49"
$ws.Range("E31").Value = "  +2.78%  "

# Row 32
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("E32").Value = "  +0.38%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.16"
$ws.Range("E33").Value = "  +8.47%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0944"
$ws.Range("E34").Value = "  +10.67%  "

# Row 35
$ws.Range("E35").Value = "  +1.56%  "

# Row 36
$ws.Range("E36").Value = "  +2.94%  "

# Row 37
$ws.Range("E37").Value = "  -0.30%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +7.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.80"
$ws.Range("E39").Value = "  -0.25%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.08"
$ws.Range("E40").Value = "  +4.89%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  +6.74%  "

# Row 42
$ws.Range("E42").Value = "  +2.08%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.60"
$ws.Range("E43").Value = "  +7.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.35"
$ws.Range("E44").Value = "  +1.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.19"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.57"
$ws.Range("E46").Value = "  +3.97%  "

# Row 47
$ws.Range("D47").Value = "2.215.45"
$ws.Range("E47").Value = "  +3.44%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("E48").Value = "  +6.73%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.270"
$ws.Range("E49").Value = "  +23.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.954"
$ws.Range("E50").Value = "  +4.07%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.58"
$ws.Range("E51").Value = "  +3.53%  "
